$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 (Viscosity) and 4:6 (Numero de Pratos, Discreto 2, Temperatura),
# leaving only the header row and the former "Densidade" row (originally row 3),
# which becomes the new row 2.
$ws.Rows("4:6").Delete()
$ws.Rows("2:2").Delete()
